# Third commit: Completed triangle.py unit tests
#
# Fill in the "Preconditions" (column E) and "Method Inputs" (column F)
# cells for the __init__ test cases of the Triangle unit-test plan:
#   - Row 7  (Attribute set to input values.)               -> Preconditions: None, Method Inputs: Shape.Triangle("red", 12, 12, 12)
#   - Row 8  (Exception raised when color is blank)          -> Preconditions: None
#   - Row 9  (Exception raised when side_1 is not an integer)-> Preconditions: None
#   - Row 10 (Exception raised when side_2 is not an integer)-> Preconditions: None
#   - Row 11 (Exception raised when side_3 is not an integer)-> Preconditions: None

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'Shape.Triangle("red", 12, 12, 12)'

$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# Leave the selection on the cell that was last edited, matching the
# author's final cursor position in the saved workbook.
$ws.Range("F7").Select() | Out-Null
